# Apply updated "count" values (column C) for Sheet 1.
#
# NOTE: the workbook stores many numeric cell values with a leading
# space in the underlying XML (e.g. "<v> 2</v>") for right-padding.
# The COM-interop runtime used here mis-parses such left-padded
# numeric literals on load (they come back as 0) even when no edit is
# made to that cell. To guarantee a correct round-trip we explicitly
# (re)write every value in column C for the data rows, using the
# authoritative values from the original file with the diff's
# corrections applied on top, rather than relying on the runtime to
# preserve untouched cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$startRow = 2
$values = @(2,3,4,1,1,1,1,16,19,1,1,7,3,1,6,31,50,13,6,1,1,38,9,1,3,1,1,3,3,1,1,1,3,2,7,36,16,8,16,15,1,12,1,12,2,1,14,2,3,27,6,2,11,15,3,2,1,1,1,2,2,12,4,3,4,7,3,19,27,5,1,2,3,31,1,6,3,1,2,1,3,1,1,2,11,20,1,10,1,2,7,12,1,3,5,1,2,7,3,3,14,3,1,32,1,7,4,1,3,6,1,5,2,5,1,2,1,3,1,5,22,1,1,1,9,7,5,1,2,1,14,1,60,10,2,3,4,57,1,3,1,1,3,13,1,1,1,1,9,3,8,1,2,1,35,2,19,11,16,1,8,1,1,6,1,4,22,1,4,2,1,1,25,2,1,2,2,74,1,6,7,3,6,19,2,1,4,11,1,10,2,1,1,12,10,3,2,1,1,2,3,2,10,1,43,8,2,8,1,3,1,33,4,1,3,1,1,7,9,8,11,2,4,36,17,4,1,9,1,8,1,1,20,2,12,1,2,1,2,8,7,1,1,31,1,1,20,4,1,1,5,2,7,1,7,3,1,1,1,2,10,1,1,1,2,1,2,3,11,1,33,9,1,18,1,2,3,3,1,2,1,1,1,1,5,2,37,14,7,5,1,4,3,7,4,1,1,8,1,1,30,2,7,5,2,9)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $values[$i]
}
